$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 to the new I1:J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-56
$iValues = @(
    7, 8, 7, 7, 6, 7, 6, 7, 6, 9, 9, 9, 7, 9, 9, 10, 8, 9, 8, 8, 7, 8, 7, 7, 8, 8, 8, 9, 7, 8, 8, 8, 8, 10, 8, 8, 7, 8, 9, 9, 8, 9, 9, 9, 8, 7, 7, 8, 6, 6, 9, 6, 3, 4, 4
)
$jValues = @(
    7, 8, 8, 8, 7, 8, 6, 8, 6, 9, 9, 9, 7, 9, 9, 10, 8, 9, 8, 8, 7, 8, 7, 7, 8, 8, 8, 9, 7, 8, 8, 8, 8, 10, 8, 8, 7, 8, 9, 9, 8, 9, 9, 9, 8, 8, 7, 8, 6, 6, 9, 6, 3, 4, 4
)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
